# Auto-generated edit script: bump "想去人数" (F column) counts
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(2, 6).Value = 39
$ws.Cells.Item(4, 6).Value = 3306
$ws.Cells.Item(5, 6).Value = 210
$ws.Cells.Item(6, 6).Value = 4808
$ws.Cells.Item(7, 6).Value = 463
$ws.Cells.Item(8, 6).Value = 291
$ws.Cells.Item(10, 6).Value = 625
$ws.Cells.Item(12, 6).Value = 25
$ws.Cells.Item(13, 6).Value = 11
$ws.Cells.Item(14, 6).Value = 651
$ws.Cells.Item(16, 6).Value = 23
$ws.Cells.Item(17, 6).Value = 89
$ws.Cells.Item(18, 6).Value = 142
$ws.Cells.Item(19, 6).Value = 343
$ws.Cells.Item(20, 6).Value = 4741
$ws.Cells.Item(24, 6).Value = 5879
$ws.Cells.Item(26, 6).Value = 1193
$ws.Cells.Item(28, 6).Value = 666
$ws.Cells.Item(30, 6).Value = 7
$ws.Cells.Item(31, 6).Value = 90
$ws.Cells.Item(32, 6).Value = 123
$ws.Cells.Item(33, 6).Value = 854
$ws.Cells.Item(34, 6).Value = 72
$ws.Cells.Item(35, 6).Value = 2
$ws.Cells.Item(36, 6).Value = 777
$ws.Cells.Item(37, 6).Value = 819

$ws = $wb.Worksheets.Item("本地生活")
$ws.Cells.Item(4, 6).Value = 37

$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(3, 6).Value = 39
$ws.Cells.Item(5, 6).Value = 37
$ws.Cells.Item(8, 6).Value = 3306
$ws.Cells.Item(9, 6).Value = 210
$ws.Cells.Item(10, 6).Value = 4808
$ws.Cells.Item(11, 6).Value = 463
$ws.Cells.Item(12, 6).Value = 291
$ws.Cells.Item(14, 6).Value = 625
$ws.Cells.Item(16, 6).Value = 25
$ws.Cells.Item(17, 6).Value = 11
$ws.Cells.Item(18, 6).Value = 651
$ws.Cells.Item(20, 6).Value = 23
$ws.Cells.Item(22, 6).Value = 89
$ws.Cells.Item(23, 6).Value = 142
$ws.Cells.Item(24, 6).Value = 343
$ws.Cells.Item(25, 6).Value = 4741
$ws.Cells.Item(29, 6).Value = 5879
$ws.Cells.Item(31, 6).Value = 1193
$ws.Cells.Item(33, 6).Value = 666
$ws.Cells.Item(35, 6).Value = 7
$ws.Cells.Item(37, 6).Value = 90
$ws.Cells.Item(38, 6).Value = 123
$ws.Cells.Item(39, 6).Value = 854
$ws.Cells.Item(40, 6).Value = 72
$ws.Cells.Item(41, 6).Value = 2
$ws.Cells.Item(42, 6).Value = 777
$ws.Cells.Item(43, 6).Value = 819
